$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as exact text without altering its style,
# avoiding Excel auto-converting numeric-looking strings to numbers.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '29.757.00'
$ws.Range("E2").Value = '  +7.05%  '
Set-TextValue $ws.Range("D3") '1.947.35'
$ws.Range("E3").Value = '  +5.38%  '
$ws.Range("E4").Value = '  -0.67%  '
Set-TextValue $ws.Range("D5") '341.52'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("E6").Value = '  -0.58%  '
Set-TextValue $ws.Range("D7") '0.4779'
$ws.Range("E7").Value = '  +2.69%  '
Set-TextValue $ws.Range("D8") '0.4129'
$ws.Range("E8").Value = '  +6.76%  '
Set-TextValue $ws.Range("D9") '47.76'
$ws.Range("E9").Value = '  +2.10%  '
Set-TextValue $ws.Range("D10") '0.08225'
$ws.Range("E10").Value = '  +3.94%  '
Set-TextValue $ws.Range("D11") '1.035'
$ws.Range("E11").Value = '  +6.74%  '
Set-TextValue $ws.Range("D12") '22.73'
$ws.Range("E12").Value = '  +6.68%  '
Set-TextValue $ws.Range("D13") '1.952.55'
$ws.Range("E13").Value = '  +5.85%  '
Set-TextValue $ws.Range("D14") '6.144'
$ws.Range("E14").Value = '  +4.10%  '
Set-TextValue $ws.Range("D15") '7.367'
$ws.Range("E15").Value = '  +2.98%  '
Set-TextValue $ws.Range("D16") '91.79'
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("E17").Value = '  -0.60%  '
Set-TextValue $ws.Range("D18") '0.00001058'
$ws.Range("E18").Value = '  +2.68%  '
Set-TextValue $ws.Range("D19") '0.06669'
$ws.Range("E19").Value = '  +0.77%  '
Set-TextValue $ws.Range("D20") '18.03'
$ws.Range("E20").Value = '  +4.27%  '
Set-TextValue $ws.Range("D22") '29.713.11'
$ws.Range("E22").Value = '  +6.85%  '
Set-TextValue $ws.Range("D23") '5.579'
$ws.Range("E23").Value = '  +4.27%  '
Set-TextValue $ws.Range("D24") '11.24'
$ws.Range("E24").Value = '  +3.69%  '
Set-TextValue $ws.Range("D25") '2.292'
$ws.Range("E25").Value = '  -0.26%  '
Set-TextValue $ws.Range("D26") '2.183.25'
$ws.Range("E26").Value = '  +5.74%  '
Set-TextValue $ws.Range("D27") '161.21'
$ws.Range("E27").Value = '  +1.46%  '
$ws.Range("E28").Value = '  +3.73%  '
Set-TextValue $ws.Range("D29") '2.167'
$ws.Range("E29").Value = '  +4.67%  '
Set-TextValue $ws.Range("D30") '5.637'
$ws.Range("E30").Value = '  +4.62%  '
Set-TextValue $ws.Range("D31") '122.99'
$ws.Range("E31").Value = '  +3.50%  '
Set-TextValue $ws.Range("D32") '1.008'
$ws.Range("E32").Value = '  +6.72%  '
Set-TextValue $ws.Range("D33") '0.09650'
$ws.Range("E33").Value = '  +2.36%  '
Set-TextValue $ws.Range("D34") '1.472'
$ws.Range("E34").Value = '  +10.67%  '
Set-TextValue $ws.Range("D35") '3.679'
$ws.Range("E35").Value = '  +2.30%  '
Set-TextValue $ws.Range("D36") '5.497'
$ws.Range("E36").Value = '  +4.44%  '
Set-TextValue $ws.Range("D37") '0.06260'
Set-TextValue $ws.Range("D38") '0.02315'
$ws.Range("E38").Value = '  +4.43%  '
Set-TextValue $ws.Range("D39") '8.494'
$ws.Range("E39").Value = '  +3.06%  '
Set-TextValue $ws.Range("D40") '1.189'
$ws.Range("E40").Value = '  +2.94%  '
Set-TextValue $ws.Range("D41") '0.6076'
$ws.Range("E41").Value = '  +4.39%  '
Set-TextValue $ws.Range("D42") '10.71'
$ws.Range("E42").Value = '  +6.27%  '
$ws.Range("E43").Value = '  -0.58%  '
Set-TextValue $ws.Range("D44") '0.1894'
$ws.Range("E44").Value = '  +2.44%  '
$ws.Range("E45").Value = '  -0.60%  '
Set-TextValue $ws.Range("D46") '2.386'
$ws.Range("E46").Value = '  +32.38%  '
Set-TextValue $ws.Range("D47") '0.5715'
$ws.Range("E47").Value = '  +4.70%  '
Set-TextValue $ws.Range("D48") '12.48'
$ws.Range("E48").Value = '  +3.95%  '
Set-TextValue $ws.Range("D49") '0.07414'
$ws.Range("E49").Value = '  +8.16%  '
Set-TextValue $ws.Range("D50") '1.990'
$ws.Range("E50").Value = '  +2.38%  '
Set-TextValue $ws.Range("D51") '112.89'
$ws.Range("E51").Value = '  +1.95%  '
